$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1738
$ws.Range("F4").Value = 124
$ws.Range("F5").Value = 348
$ws.Range("F6").Value = 745
$ws.Range("F7").Value = 193
$ws.Range("F8").Value = 239
$ws.Range("F9").Value = 1034
$ws.Range("F12").Value = 600
$ws.Range("F14").Value = 490
$ws.Range("F17").Value = 152
$ws.Range("F18").Value = 808
$ws.Range("F19").Value = 2588
$ws.Range("F20").Value = 520
$ws.Range("F24").Value = 198
$ws.Range("F26").Value = 144
$ws.Range("F29").Value = 54
$ws.Range("F30").Value = 198
$ws.Range("F31").Value = 1011
$ws.Range("F32").Value = 69
$ws.Range("F33").Value = 39

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1017
$ws.Range("F5").Value = 1017
$ws.Range("F10").Value = 308
$ws.Range("F14").Value = 563
$ws.Range("F15").Value = 88
$ws.Range("F17").Value = 963
$ws.Range("F20").Value = 607
$ws.Range("F21").Value = 2
$ws.Range("F24").Value = 277
$ws.Range("F25").Value = 244
$ws.Range("F26").Value = 1543
$ws.Range("F29").Value = 16
$ws.Range("F31").Value = 4
$ws.Range("F33").Value = 90

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1765
$ws.Range("F5").Value = 2379
$ws.Range("F6").Value = 970
$ws.Range("F9").Value = 1221
$ws.Range("F11").Value = 87

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1765
$ws.Range("F5").Value = 2379
$ws.Range("F7").Value = 1738
$ws.Range("F8").Value = 970
$ws.Range("F9").Value = 1221
$ws.Range("F11").Value = 87
$ws.Range("F12").Value = 124
$ws.Range("F13").Value = 348
$ws.Range("F14").Value = 745
$ws.Range("F15").Value = 193
$ws.Range("F17").Value = 239
$ws.Range("F18").Value = 1034
$ws.Range("F20").Value = 600
$ws.Range("F21").Value = 1017
$ws.Range("F22").Value = 490
$ws.Range("F23").Value = 134
$ws.Range("F25").Value = 152
$ws.Range("F26").Value = 808
$ws.Range("F27").Value = 2588
$ws.Range("F28").Value = 520
$ws.Range("F31").Value = 198
$ws.Range("F32").Value = 144
$ws.Range("F35").Value = 563
$ws.Range("F36").Value = 88
$ws.Range("F37").Value = 54
$ws.Range("F38").Value = 198
$ws.Range("F41").Value = 2
$ws.Range("F43").Value = 277
$ws.Range("F44").Value = 277
$ws.Range("F45").Value = 244
$ws.Range("F46").Value = 1011
